$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell F6 text
$ws.Range("F6").Value = "master_all_responses_click_worker_44HITsOct-01-2023.csv"

# Update selection to F9
$ws.Range("F9").Select()

# Update column F width, remove autofit/bestfit
# Note: the runtime's ColumnWidth setter applies Excel's standard +5px
# padding (at 6px/char MDW here) before storing the XML "width" attribute,
# so feed it a value that round-trips to exactly 75.5 after that padding.
$ws.Columns("F").ColumnWidth = 74.66666666666667
